# Apply the "Modify undo redo diagrams" edit:
#  1. Update the cached date placeholder text on every slide layout
#     from "6/7/2018" to "14/4/19".
#  2. Rename "address book" -> "card collection" (and the matching
#     camelCase identifier) in the two activity-diagram labels on
#     slide 1.

$p = $ppt.ActivePresentation

# --- 1. Slide layout date placeholders -------------------------------
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
        }
        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "6/7/2018") {
                $tr.Text = "14/4/19"
            }
        }
    }
}

# --- 2. Slide 1 wording updates ---------------------------------------
# Runs are replaced whole (matching their original run boundaries) so
# that the underlying <a:r> split is preserved exactly, instead of
# fragmenting a run by overwriting only part of its text.
$slide = $p.Slides.Item(1)

# "[command commits address book]" -> "[command commits card collection]"
# Shape text is: "[" (run 1) + "command commits address book]" (run 2)
# The shape uses spAutoFit, so editing the text makes the host
# recompute the box height; put the original height back afterwards so
# the only change is the text itself, same as the source edit.
$shGuardText = $slide.Shapes.Item("TextBox 47")
$origHeight = $shGuardText.Height

$trGuard = $shGuardText.TextFrame.TextRange
$run2Guard = $trGuard.Characters(2, $trGuard.Length - 1)
$run2Guard.Text = "command commits card collection]"

$shGuardText.Height = $origHeight

# "Purge redundant states and then save address book to " (run 1)
#  + "addressBookStateList" (run 2, err="1") + " " (run 3)
# ->
# "Purge redundant states and then save card collection to " (run 1)
#  + "cardCollectionStateList" (run 2) + " " (run 3, untouched)
$shPurge = $slide.Shapes.Item("Rounded Rectangle 50")
$trPurge = $shPurge.TextFrame.TextRange
$newRun1Text = "Purge redundant states and then save card collection to "
$run1Purge = $trPurge.Characters(1, 53)
$run1Purge.Text = $newRun1Text

$trPurge2 = $shPurge.TextFrame.TextRange
$run2Purge = $trPurge2.Characters($newRun1Text.Length + 1, 20)
$run2Purge.Text = "cardCollectionStateList"
